$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.013.93'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").Value = '2.624.08'
$ws.Range("E3").Value = '  -1.40%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.584'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").Value = '2.624.89'
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.361'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").Value = '3.094.68'
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '62.909.21'
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("E17").Value = '  -2.04%  '
$ws.Range("D18").Value = '2.637.06'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("E24").Value = '  -4.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("B28").Value = 'SuiNetwork'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.89%  '
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '536.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.83'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("D35").Value = '0.0₃0799'
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("E36").Value = '  +12.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.401'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("E41").Value = '  +6.33%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '167.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.67%  '
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0564'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.623'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.31%  '
$ws.Range("E51").Value = '  -0.34%  '
